$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.190.14'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.654.29'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.38%  '
$rng = $ws.Range('D5')
$rng.NumberFormat = "@"
$rng.Value = '218.12'
$rng.Style = "Normal"
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('E7').Value = '  +0.35%  '
$rng = $ws.Range('D8')
$rng.NumberFormat = "@"
$rng.Value = '0.2646'
$rng.Style = "Normal"
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('E9').Value = '  -1.25%  '
$rng = $ws.Range('D10')
$rng.NumberFormat = "@"
$rng.Value = '21.24'
$rng.Style = "Normal"
$ws.Range('E10').Value = '  +1.65%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '1.656.35'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E13').Value = '  -0.95%  '
$rng = $ws.Range('D14')
$rng.NumberFormat = "@"
$rng.Value = '0.5458'
$rng.Style = "Normal"
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Value = '0.0₅8210'
$ws.Range('E15').Value = '  -0.76%  '
$rng = $ws.Range('D16')
$rng.NumberFormat = "@"
$rng.Value = '64.71'
$rng.Style = "Normal"
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').Value = '26.202.50'
$rng = $ws.Range('D19')
$rng.NumberFormat = "@"
$rng.Value = '4.679'
$rng.Style = "Normal"
$ws.Range('E19').Value = '  -1.69%  '
$rng = $ws.Range('D20')
$rng.NumberFormat = "@"
$rng.Value = '191.44'
$rng.Style = "Normal"
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('E21').Value = '  -1.40%  '
$rng = $ws.Range('D22')
$rng.NumberFormat = "@"
$rng.Value = '6.181'
$rng.Style = "Normal"
$ws.Range('E22').Value = '  -2.86%  '
$rng = $ws.Range('D23')
$rng.NumberFormat = "@"
$rng.Value = '1.007'
$rng.Style = "Normal"
$ws.Range('E23').Value = '  +0.50%  '
$rng = $ws.Range('D24')
$rng.NumberFormat = "@"
$rng.Value = '138.47'
$rng.Style = "Normal"
$ws.Range('E24').Value = '  -3.26%  '
$rng = $ws.Range('D25')
$rng.NumberFormat = "@"
$rng.Value = '0.1241'
$rng.Style = "Normal"
$ws.Range('E25').Value = '  -1.05%  '
$rng = $ws.Range('D26')
$rng.NumberFormat = "@"
$rng.Value = '7.283'
$rng.Style = "Normal"
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('E27').Value = '  +0.53%  '
$rng = $ws.Range('D28')
$rng.NumberFormat = "@"
$rng.Value = '1.415'
$rng.Style = "Normal"
$ws.Range('E28').Value = '  -0.97%  '
$rng = $ws.Range('D29')
$rng.NumberFormat = "@"
$rng.Value = '0.06067'
$rng.Style = "Normal"
$ws.Range('E29').Value = '  -1.24%  '
$rng = $ws.Range('D30')
$rng.NumberFormat = "@"
$rng.Value = '1.283'
$rng.Style = "Normal"
$ws.Range('E30').Value = '  +1.35%  '
$rng = $ws.Range('D31')
$rng.NumberFormat = "@"
$rng.Value = '3.540'
$rng.Style = "Normal"
$ws.Range('E31').Value = '  -0.43%  '
$rng = $ws.Range('D32')
$rng.NumberFormat = "@"
$rng.Value = '3.360'
$rng.Style = "Normal"
$ws.Range('E32').Value = '  -1.85%  '
$rng = $ws.Range('D33')
$rng.NumberFormat = "@"
$rng.Value = '1.653'
$rng.Style = "Normal"
$ws.Range('E33').Value = '  -0.59%  '
$rng = $ws.Range('D34')
$rng.NumberFormat = "@"
$rng.Value = '0.9845'
$rng.Style = "Normal"
$ws.Range('E34').Value = '  -1.51%  '
$rng = $ws.Range('D35')
$rng.NumberFormat = "@"
$rng.Value = '2.411'
$rng.Style = "Normal"
$ws.Range('E35').Value = '  +0.45%  '
$rng = $ws.Range('D36')
$rng.NumberFormat = "@"
$rng.Value = '2.769'
$rng.Style = "Normal"
$ws.Range('E36').Value = '  +0.25%  '
$rng = $ws.Range('D37')
$rng.NumberFormat = "@"
$rng.Value = '0.5943'
$rng.Style = "Normal"
$ws.Range('E37').Value = '  +4.83%  '
$ws.Range('E38').Value = '  -0.58%  '
$rng = $ws.Range('D39')
$rng.NumberFormat = "@"
$rng.Value = '5.961'
$rng.Style = "Normal"
$ws.Range('E39').Value = '  +0.78%  '
$rng = $ws.Range('D40')
$rng.NumberFormat = "@"
$rng.Value = '0.8627'
$rng.Style = "Normal"
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('D41').Value = '1.058.06'
$ws.Range('E41').Value = '  +2.53%  '
$rng = $ws.Range('D42')
$rng.NumberFormat = "@"
$rng.Value = '1.003'
$rng.Style = "Normal"
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '1.796.87'
$ws.Range('E44').Value = '  -0.46%  '
$rng = $ws.Range('D45')
$rng.NumberFormat = "@"
$rng.Value = '57.39'
$rng.Style = "Normal"
$ws.Range('E45').Value = '  +2.26%  '
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('E47').Value = '  +0.04%  '
$rng = $ws.Range('D48')
$rng.NumberFormat = "@"
$rng.Value = '8.072'
$rng.Style = "Normal"
$ws.Range('E48').Value = '  -0.48%  '
$rng = $ws.Range('D49')
$rng.NumberFormat = "@"
$rng.Value = '0.05179'
$rng.Style = "Normal"
$ws.Range('E49').Value = '  +0.28%  '
$rng = $ws.Range('D50')
$rng.NumberFormat = "@"
$rng.Value = '1.465'
$rng.Style = "Normal"
$ws.Range('E50').Value = '  +4.90%  '
$rng = $ws.Range('D51')
$rng.NumberFormat = "@"
$rng.Value = '0.4231'
$rng.Style = "Normal"
$ws.Range('E51').Value = '  +0.38%  '
